# Updated cryptos list data (price + 1h volume change), scraped values refreshed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.879.33"
$ws.Range("E2").Value = "  +3.48%  "
$ws.Range("D3").Value = "2.621.96"
$ws.Range("E3").Value = "  +4.65%  "
$ws.Range("E4").Value = "  -0.21%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "326.95"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.92%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "109.87"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.95%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.533"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.07%  "
$ws.Range("E8").Value = "  -0.11%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.559"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.60%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "40.10"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.81%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "20.60"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("E12").Value = "  +0.31%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.26"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").Value = "3.032.18"
$ws.Range("E15").Value = "  +4.41%  "
$ws.Range("D16").Value = "2.617.99"
$ws.Range("E16").Value = "  +4.17%  "
$ws.Range("E17").Value = "  +3.71%  "
$ws.Range("D18").Value = "49.775.75"
$ws.Range("E18").Value = "  +3.48%  "
$ws.Range("E19").Value = "  +11.46%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "13.32"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.66%  "
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("E22").Value = "  +1.72%  "
$ws.Range("E23").Value = "  +0.81%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "278.06"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("E26").Value = "  +3.09%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  +0.77%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.95"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.86%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "36.56"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +3.84%  "
$ws.Range("E31").Value = "  +1.97%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "49.77"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.56%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "19.70"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.10%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.43"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("E35").Value = "  -0.16%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.0789"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("E37").Value = "  +5.45%  "
$ws.Range("E38").Value = "  +1.99%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "3.13"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +7.27%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "123.88"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("E41").Value = "  +0.74%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "22.53"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +4.18%  "
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("E44").Value = "  +4.09%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "3.35"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +5.34%  "
$ws.Range("D46").Value = "2.048.41"
$ws.Range("E46").Value = "  +1.66%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.31"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +15.84%  "
$ws.Range("E48").Value = "  +8.65%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "9.04"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.28%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "5.35"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.29%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "81.50"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.21%  "
